$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the hotel id "10001" onto each existing id_servicios value in column A
$ws.Range("A2").Value = 1000110101
$ws.Range("A3").Value = 1000110102
$ws.Range("A4").Value = 1000110103
$ws.Range("A5").Value = 1000110104
$ws.Range("A6").Value = 1000110105
$ws.Range("A7").Value = 1000110106
$ws.Range("A8").Value = 1000110107
$ws.Range("A9").Value = 1000110108
$ws.Range("A10").Value = 1000110109
$ws.Range("A11").Value = 1000110110

# Widen column A to fit the longer ids
$ws.Columns.Item(1).ColumnWidth = 11.1640625

# Move the active selection to B14 (matches the author's last cursor position)
$ws.Range("B14").Select()
